# Update crypto price/volume figures per the Sep 5 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column values look like plain decimals (e.g. "215.24"); Excel's
# normal text-entry auto-detection would turn those into numbers and silently
# drop significant trailing zeros (e.g. "6.840" -> 6.84). Force those specific
# cells to Text format first so the literal string is preserved, exactly as
# the source data (and the original workbook) stores them.
$textForceCells = @(
    "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D15", "D17", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell -> new value map (applied in sheet order).
$ws.Range("D2").Value = "25.835.69"
$ws.Range("D3").Value = "1.636.06"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "215.24"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "0.5086"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").Value = "0.06428"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").Value = "20.39"
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("D11").Value = "0.07806"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "1.655.71"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").Value = "4.259"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "1.863.92"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "0.5594"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "0.0₅7664"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "63.27"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "25.861.23"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "193.31"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "4.382"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").Value = "6.148"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "1.757"
$ws.Range("E25").Value = "  -6.89%  "
$ws.Range("D26").Value = "138.47"
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").Value = "6.840"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").Value = "15.57"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "0.04967"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "3.301"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").Value = "3.254"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("D34").Value = "1.570"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("D35").Value = "2.389"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "0.9028"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").Value = "2.578"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("D38").Value = "0.5561"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").Value = "1.133.59"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("D40").Value = "0.01571"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("D41").Value = "0.9968"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "5.463"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "99.12"
$ws.Range("D44").Value = "0.8007"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "55.64"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").Value = "0.4269"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("D48").Value = "7.788"
$ws.Range("E48").Value = "  +3.82%  "
$ws.Range("D49").Value = "0.05031"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  +0.48%  "
